$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.807599666666667
$ws.Range("H2").Value = 5.422799
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 14.61878266666667
$ws.Range("N2").Value = 43.856348
$ws.Range("O2").Value = 0.2662829816142094
$ws.Range("P2").Value = 0.2662829816142094
$ws.Range("Q2").Value = 26.42490667533911
$ws.Range("R2").Value = 237.824160078052
$ws.Range("S2").Value = 0.2662829816142094
$ws.Range("T2").Value = 0.2662829816142094

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.807599666666667
$ws.Range("H3").Value = 5.422799
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 27.084169
$ws.Range("N3").Value = 81.25250700000001
$ws.Range("O3").Value = 0.4933415757187404
$ws.Range("P3").Value = 0.4933415757187404
$ws.Range("Q3").Value = 48.95733485634367
$ws.Range("R3").Value = 440.6160137070931
$ws.Range("S3").Value = 0.4933415757187404
$ws.Range("T3").Value = 0.4933415757187404

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.807599666666667
$ws.Range("H4").Value = 5.422799
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 13.19647366666667
$ws.Range("N4").Value = 39.589421
$ws.Range("O4").Value = 0.2403754426670501
$ws.Range("P4").Value = 0.2403754426670501
$ws.Range("Q4").Value = 23.85394140104211
$ws.Range("R4").Value = 214.685472609379
$ws.Range("S4").Value = 0.2403754426670501
$ws.Range("T4").Value = 0.2403754426670501
